$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Extend header row (row 1) with the bold/bordered header style ---
$ws.Range("B1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)

# --- Extend data row (row 2) with the data-row style ---
$ws.Range("B2").Copy()
$ws.Range("E2:K2").PasteSpecial(-4122)

# --- Fix up the existing header cells (they previously held stray data values) ---
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"

# --- Fill in the new header cells ---
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Fill in the new data cells for row 2 ---
$ws.Range("E2").Value = "insurance"
$ws.Range("F2").Value = "normal"
$ws.Range("G2").Value = "2013-12-31"
$ws.Range("H2").Value = "許添財"
$ws.Range("I2").Value = 639
$ws.Range("J2").Value = "tmpbb0f1"
$ws.Range("K2").Value = 94

Write-Host "done"
